$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1516543333333333
$ws.Range("H2").Value = 0.454963
$ws.Range("I2").Value = 0.0007044400935133411
$ws.Range("J2").Value = 0.0007044400935133412
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03186166666666667
$ws.Range("N2").Value = 0.095585
$ws.Range("O2").Value = 0.0002078156820111728
$ws.Range("P2").Value = 0.0002078156820111728
$ws.Range("Q2").Value = 0.004831959817222222
$ws.Range("R2").Value = 0.043487638355
$ws.Range("S2").Value = 0.0000001463936984694893
$ws.Range("T2").Value = 0.0000001463936984694893
$ws.Range("G3").Value = 0.1516543333333333
$ws.Range("H3").Value = 0.454963
$ws.Range("I3").Value = 0.0007044400935133411
$ws.Range("J3").Value = 0.0007044400935133412
$ws.Range("O3").Value = 0.0004817036408055181
$ws.Range("P3").Value = 0.0004817036408055181
$ws.Range("Q3").Value = 0.01120017803111111
$ws.Range("R3").Value = 0.10080160228
$ws.Range("S3").Value = 0.0000003393313577747561
$ws.Range("T3").Value = 0.0000003393313577747561
$ws.Range("G4").Value = 0.1516543333333333
$ws.Range("H4").Value = 0.454963
$ws.Range("I4").Value = 0.0007044400935133411
$ws.Range("J4").Value = 0.0007044400935133412
$ws.Range("M4").Value = 0.1279203333333333
$ws.Range("N4").Value = 0.383761
$ws.Range("O4").Value = 0.0008343521885681821
$ws.Range("P4").Value = 0.000834352188568182
$ws.Range("Q4").Value = 0.01939967287144444
$ws.Range("R4").Value = 0.174597055843
$ws.Range("S4").Value = 0.000000587751133738031
$ws.Range("T4").Value = 0.0000005877511337380311
$ws.Range("G5").Value = 0.1516543333333333
$ws.Range("H5").Value = 0.454963
$ws.Range("I5").Value = 0.0007044400935133411
$ws.Range("J5").Value = 0.0007044400935133412
$ws.Range("M5").Value = 153.0833153333333
$ws.Range("N5").Value = 459.249946
$ws.Range("O5").Value = 0.9984761284886152
$ws.Range("P5").Value = 0.9984761284886152
$ws.Range("Q5").Value = 23.21574813133311
$ws.Range("R5").Value = 208.941733181998
$ws.Range("S5").Value = 0.000703366617323359
$ws.Range("T5").Value = 0.0007033666173233591
$ws.Range("I6").Value = 0.0008677905854558892
$ws.Range("J6").Value = 0.0008677905854558892
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03186166666666667
$ws.Range("N6").Value = 0.095585
$ws.Range("O6").Value = 0.0002078156820111728
$ws.Range("P6").Value = 0.0002078156820111728
$ws.Range("Q6").Value = 0.005952428428333333
$ws.Range("R6").Value = 0.053571855855
$ws.Range("S6").Value = 0.0000001803404923593906
$ws.Range("T6").Value = 0.0000001803404923593905
$ws.Range("I7").Value = 0.0008677905854558892
$ws.Range("J7").Value = 0.0008677905854558892
$ws.Range("O7").Value = 0.0004817036408055181
$ws.Range("P7").Value = 0.0004817036408055181
$ws.Range("S7").Value = 0.0000004180178844708539
$ws.Range("T7").Value = 0.0000004180178844708539
$ws.Range("I8").Value = 0.0008677905854558892
$ws.Range("J8").Value = 0.0008677905854558892
$ws.Range("M8").Value = 0.1279203333333333
$ws.Range("N8").Value = 0.383761
$ws.Range("O8").Value = 0.0008343521885681821
$ws.Range("P8").Value = 0.000834352188568182
$ws.Range("Q8").Value = 0.02389820459366666
$ws.Range("R8").Value = 0.215083841343
$ws.Range("S8").Value = 0.0000007240429741939853
$ws.Range("T8").Value = 0.0000007240429741939852
$ws.Range("I9").Value = 0.0008677905854558892
$ws.Range("J9").Value = 0.0008677905854558892
$ws.Range("M9").Value = 153.0833153333333
$ws.Range("N9").Value = 459.249946
$ws.Range("O9").Value = 0.9984761284886152
$ws.Range("P9").Value = 0.9984761284886152
$ws.Range("Q9").Value = 28.59917805388866
$ws.Range("R9").Value = 257.392602484998
$ws.Range("S9").Value = 0.000866468184104865
$ws.Range("T9").Value = 0.000866468184104865
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.003907
$ws.Range("H10").Value = 0.011721
$ws.Range("I10").Value = 0.00001814816223752233
$ws.Range("J10").Value = 0.00001814816223752233
$ws.Range("M10").Value = 0.03186166666666667
$ws.Range("N10").Value = 0.095585
$ws.Range("O10").Value = 0.0002078156820111728
$ws.Range("P10").Value = 0.0002078156820111728
$ws.Range("Q10").Value = 0.0001244835316666667
$ws.Range("R10").Value = 0.001120351785
$ws.Range("S10").Value = 0.000000003771472712640115
$ws.Range("T10").Value = 0.000000003771472712640115
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.003907
$ws.Range("H11").Value = 0.011721
$ws.Range("I11").Value = 0.00001814816223752233
$ws.Range("J11").Value = 0.00001814816223752233
$ws.Range("O11").Value = 0.0004817036408055181
$ws.Range("P11").Value = 0.0004817036408055181
$ws.Range("Q11").Value = 0.0002885449733333334
$ws.Range("R11").Value = 0.00259690476
$ws.Range("S11").Value = 0.000000008742035823743724
$ws.Range("T11").Value = 0.000000008742035823743724
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.003907
$ws.Range("H12").Value = 0.011721
$ws.Range("I12").Value = 0.00001814816223752233
$ws.Range("J12").Value = 0.00001814816223752233
$ws.Range("M12").Value = 0.1279203333333333
$ws.Range("N12").Value = 0.383761
$ws.Range("O12").Value = 0.0008343521885681821
$ws.Range("P12").Value = 0.000834352188568182
$ws.Range("Q12").Value = 0.0004997847423333333
$ws.Range("R12").Value = 0.004498062681
$ws.Range("S12").Value = 0.00000001514195888136719
$ws.Range("T12").Value = 0.00000001514195888136719
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.003907
$ws.Range("H13").Value = 0.011721
$ws.Range("I13").Value = 0.00001814816223752233
$ws.Range("J13").Value = 0.00001814816223752233
$ws.Range("M13").Value = 153.0833153333333
$ws.Range("N13").Value = 459.249946
$ws.Range("O13").Value = 0.9984761284886152
$ws.Range("P13").Value = 0.9984761284886152
$ws.Range("Q13").Value = 0.5980965130073334
$ws.Range("R13").Value = 5.382868617066
$ws.Range("S13").Value = 0.00001812050677010458
$ws.Range("T13").Value = 0.00001812050677010458
$ws.Range("G14").Value = 214.941124
$ws.Range("H14").Value = 644.8233720000001
$ws.Range("I14").Value = 0.9984096211587932
$ws.Range("J14").Value = 0.9984096211587933
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.03186166666666667
$ws.Range("N14").Value = 0.095585
$ws.Range("O14").Value = 0.0002078156820111728
$ws.Range("P14").Value = 0.0002078156820111728
$ws.Range("Q14").Value = 6.848382445846668
$ws.Range("R14").Value = 61.63544201262001
$ws.Range("S14").Value = 0.0002074851763476313
$ws.Range("T14").Value = 0.0002074851763476313
$ws.Range("G15").Value = 214.941124
$ws.Range("H15").Value = 644.8233720000001
$ws.Range("I15").Value = 0.9984096211587932
$ws.Range("J15").Value = 0.9984096211587933
$ws.Range("O15").Value = 0.0004817036408055181
$ws.Range("P15").Value = 0.0004817036408055181
$ws.Range("Q15").Value = 15.87411847781334
$ws.Range("R15").Value = 142.86706630032
$ws.Range("S15").Value = 0.0004809375495274487
$ws.Range("T15").Value = 0.0004809375495274487
$ws.Range("G16").Value = 214.941124
$ws.Range("H16").Value = 644.8233720000001
$ws.Range("I16").Value = 0.9984096211587932
$ws.Range("J16").Value = 0.9984096211587933
$ws.Range("M16").Value = 0.1279203333333333
$ws.Range("N16").Value = 0.383761
$ws.Range("O16").Value = 0.0008343521885681821
$ws.Range("P16").Value = 0.000834352188568182
$ws.Range("Q16").Value = 27.49534022912134
$ws.Range("R16").Value = 247.458062062092
$ws.Range("S16").Value = 0.0008330252525013687
$ws.Range("T16").Value = 0.0008330252525013687
$ws.Range("G17").Value = 214.941124
$ws.Range("H17").Value = 644.8233720000001
$ws.Range("I17").Value = 0.9984096211587932
$ws.Range("J17").Value = 0.9984096211587933
$ws.Range("M17").Value = 153.0833153333333
$ws.Range("N17").Value = 459.249946
$ws.Range("O17").Value = 0.9984761284886152
$ws.Range("P17").Value = 0.9984761284886152
$ws.Range("Q17").Value = 32903.8998633931
$ws.Range("R17").Value = 296135.0987705379
$ws.Range("S17").Value = 0.9968881731804168
$ws.Range("T17").Value = 0.9968881731804169
